# Regional Availability Factor workbook update:
# Split the single "RAF" sheet into "RAF-generation" (existing data) and a
# new "RAF-demand-altering-techs" sheet, and update the About sheet titles
# to reflect the two RAF variables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing "RAF" sheet to "RAF-generation"
# ---------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("RAF")
$wsGen.Name = "RAF-generation"

# Reset its selection back to the default A1 cell (was B2:B24)
[void]$wsGen.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Add the new "RAF-demand-altering-techs" sheet right after it
# ---------------------------------------------------------------------
$wsDemand = $wb.Worksheets.Add($null, $wsGen)
$wsDemand.Name = "RAF-demand-altering-techs"

# Match the look of the generation sheet (tab color + column widths)
$wsDemand.Tab.Color = $wsGen.Tab.Color
$wsDemand.Columns.Item(1).ColumnWidth = $wsGen.Columns.Item(1).ColumnWidth
$wsDemand.Columns.Item(2).ColumnWidth = 19.140625

# Header row
$wsDemand.Range("A1").Value = "Unit: dimensionless (% of capacity available)"
$wsDemand.Range("B1").Value = "Percent of capacity"
$wsDemand.Range("A1").Font.Italic = $true
$wsDemand.Range("A1").WrapText = $true
$wsDemand.Rows.Item(1).RowHeight = 30

# Data row
$wsDemand.Range("A2").Value = "demand-altering technologies"
$wsDemand.Range("B2").Value = 0.9

[void]$wsDemand.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Update the "About" sheet: split the single title into two titles,
#    one for each RAF sheet
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A1").Value = "RAF Regional Availability Factor for Generation"
[void]$wsAbout.Rows.Item(2).Insert()
$wsAbout.Range("A2").Value = "RAF Regional Availability Factor for Demand-Altering Technologies"

[void]$wsAbout.Range("A1").Select()
[void]$wsAbout.Activate()

Write-Host "RAF workbook split into generation / demand-altering-techs sheets"
